$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.71700912156666
$ws.Range("C2").Value = 6.02658000870023
$ws.Range("D2").Value = 7.079662815943526
$ws.Range("F2").Value = 43.76194602766287
$ws.Range("G2").Value = 52.79316675476553
$ws.Range("H2").Value = 20.37363761606754
$ws.Range("I2").Value = 32.62031354980262
$ws.Range("J2").Value = 10.80047491848372
$ws.Range("K2").Value = 11.32208869761807
$ws.Range("L2").Value = 11.79266655893097
$ws.Range("M2").Value = 15.9620900634532
$ws.Range("B3").Value = 13.57328223772181
$ws.Range("C3").Value = 5.927241403103127
$ws.Range("D3").Value = 7.076893557270616
$ws.Range("F3").Value = 43.82139688016896
$ws.Range("G3").Value = 52.84427887274726
$ws.Range("H3").Value = 20.4175095719532
$ws.Range("I3").Value = 32.69445577979634
$ws.Range("J3").Value = 10.81735564033697
$ws.Range("K3").Value = 11.22060214292735
$ws.Range("L3").Value = 11.80503006688283
$ws.Range("M3").Value = 15.95377303742342
$ws.Range("B4").Value = 13.48719281285092
$ws.Range("C4").Value = 5.864399611568938
$ws.Range("D4").Value = 7.075978155568792
$ws.Range("F4").Value = 43.8656968381714
$ws.Range("G4").Value = 52.88670086973527
$ws.Range("H4").Value = 20.4472359031049
$ws.Range("I4").Value = 32.74482629668985
$ws.Range("J4").Value = 10.82831525157785
$ws.Range("K4").Value = 11.16015640888212
$ws.Range("L4").Value = 11.81396791122716
$ws.Range("M4").Value = 15.95093153853174
$ws.Range("B5").Value = 13.45269233454891
$ws.Range("C5").Value = 5.838340374768144
$ws.Range("D5").Value = 7.075803524645825
$ws.Range("F5").Value = 43.88570844520562
$ws.Range("G5").Value = 52.90676037262596
$ws.Range("H5").Value = 20.46005077024401
$ws.Range("I5").Value = 32.76657071043153
$ws.Range("J5").Value = 10.8329313524803
$ws.Range("K5").Value = 11.13601897014913
$ws.Range("L5").Value = 11.81794934142171
$ws.Range("M5").Value = 15.95034571227713
$ws.Range("B6").Value = 13.44699979550475
$ws.Range("C6").Value = 5.833986454052371
$ws.Range("D6").Value = 7.075786539184273
$ws.Range("F6").Value = 43.88914962226826
$ws.Range("G6").Value = 52.91025854577293
$ws.Range("H6").Value = 20.46222101273931
$ws.Range("I6").Value = 32.77025489533095
$ws.Range("J6").Value = 10.83370692237462
$ws.Range("K6").Value = 11.13204156099413
$ws.Range("L6").Value = 11.81863095691985
$ws.Range("M6").Value = 15.95028305893264
$ws.Range("B7").Value = 13.48672512195119
$ws.Range("C7").Value = 5.864049973685511
$ws.Range("D7").Value = 7.075974995789383
$ws.Range("F7").Value = 43.86595879217602
$ws.Range("G7").Value = 52.88696018029898
$ws.Range("H7").Value = 20.44740589016455
$ws.Range("I7").Value = 32.74511461884042
$ws.Range("J7").Value = 10.82837689816162
$ws.Range("K7").Value = 11.15982884745996
$ws.Range("L7").Value = 11.81402023210268
$ws.Range("M7").Value = 15.9509213182906
$ws.Range("B8").Value = 13.66702739306877
$ws.Range("C8").Value = 5.992719635041855
$ws.Range("D8").Value = 7.078545692308008
$ws.Range("F8").Value = 43.78082560887087
$ws.Range("G8").Value = 52.80849717316185
$ws.Range("H8").Value = 20.38818579707854
$ws.Range("I8").Value = 32.64487137516189
$ws.Range("J8").Value = 10.80617219073223
$ws.Range("K8").Value = 11.28672387824642
$ws.Range("L8").Value = 11.79665036227665
$ws.Range("M8").Value = 15.95875378150918
$ws.Range("B9").Value = 14.03599091496534
$ws.Range("C9").Value = 6.229805180634774
$ws.Range("D9").Value = 7.089768353275164
$ws.Range("F9").Value = 43.67580125424789
$ws.Range("G9").Value = 52.74234641380402
$ws.Range("H9").Value = 20.29418770939654
$ws.Range("I9").Value = 32.48678564418521
$ws.Range("J9").Value = 10.76732958414804
$ws.Range("K9").Value = 11.54923965972021
$ws.Range("L9").Value = 11.77324699308878
$ws.Range("M9").Value = 15.99196847133803
$ws.Range("B10").Value = 14.31407238711708
$ws.Range("C10").Value = 6.393980277979887
$ws.Range("D10").Value = 7.101717285377858
$ws.Range("F10").Value = 43.63644565847405
$ws.Range("G10").Value = 52.74733614541102
$ws.Range("H10").Value = 20.23862548203823
$ws.Range("I10").Value = 32.39414831279491
$ws.Range("J10").Value = 10.74163175916797
$ws.Range("K10").Value = 11.74887290273452
$ws.Range("L10").Value = 11.7625148734282
$ws.Range("M10").Value = 16.02708654865937
$ws.Range("B11").Value = 14.44158858992703
$ws.Range("C11").Value = 6.466353713448672
$ws.Range("D11").Value = 7.107942414179172
$ws.Range("F11").Value = 43.62675347501307
$ws.Range("G11").Value = 52.76124579001613
$ws.Range("H11").Value = 20.216280415216
$ws.Range("I11").Value = 32.3571176781072
$ws.Range("J11").Value = 10.73055231106809
$ws.Range("K11").Value = 11.8408166086522
$ws.Range("L11").Value = 11.75902746452257
$ws.Range("M11").Value = 16.04534845089827
$ws.Range("B12").Value = 14.4899790567539
$ws.Range("C12").Value = 6.493416451237246
$ws.Range("D12").Value = 7.11041185458329
$ws.Range("F12").Value = 43.62426322137135
$ws.Range("G12").Value = 52.76818491063668
$ws.Range("H12").Value = 20.20824031408998
$ws.Range("I12").Value = 32.34383055390912
$ws.Range("J12").Value = 10.726444203857
$ws.Range("K12").Value = 11.87576670478954
$ws.Range("L12").Value = 11.75790664297966
$ws.Range("M12").Value = 16.05258865820717
$ws.Range("B13").Value = 14.47955341197451
$ws.Range("C13").Value = 6.487603456011312
$ws.Range("D13").Value = 7.109875052965052
$ws.Range("F13").Value = 43.6247470804173
$ws.Range("G13").Value = 52.76661613495718
$ws.Range("H13").Value = 20.20995314678915
$ws.Range("I13").Value = 32.34665944913102
$ws.Range("J13").Value = 10.72732507517351
$ws.Range("K13").Value = 11.8682341280809
$ws.Range("L13").Value = 11.75813915921655
$ws.Range("M13").Value = 16.0510149693402
$ws.Range("B14").Value = 14.44556788712499
$ws.Range("C14").Value = 6.468587127643485
$ws.Range("D14").Value = 7.108143337779143
$ws.Range("F14").Value = 43.62652495828365
$ws.Range("G14").Value = 52.76178318485913
$ws.Range("H14").Value = 20.21561050357861
$ws.Range("I14").Value = 32.3560097928366
$ws.Range("J14").Value = 10.73021258429281
$ws.Range("K14").Value = 11.8436894765896
$ws.Range("L14").Value = 11.7589312544481
$ws.Range("M14").Value = 16.04593762082696
$ws.Range("B15").Value = 14.42476289231675
$ws.Range("C15").Value = 6.456894019993984
$ws.Range("D15").Value = 7.107097170592942
$ws.Range("F15").Value = 43.62776759476195
$ws.Range("G15").Value = 52.75904050823471
$ws.Range("H15").Value = 20.21913069051367
$ws.Range("I15").Value = 32.36183296039081
$ws.Range("J15").Value = 10.73199264334269
$ws.Range("K15").Value = 11.82867161886694
$ws.Range("L15").Value = 11.75944242988357
$ws.Range("M15").Value = 16.04286977751418
$ws.Range("B16").Value = 14.30575563675746
$ws.Range("C16").Value = 6.389203031056728
$ws.Range("D16").Value = 7.101326226180207
$ws.Range("F16").Value = 43.63724423764813
$ws.Range("G16").Value = 52.74666120035332
$ws.Range("H16").Value = 20.24014476215424
$ws.Range("I16").Value = 32.39667124965954
$ws.Range("J16").Value = 10.74236808294716
$ws.Range("K16").Value = 11.74288437885262
$ws.Range("L16").Value = 11.76277078368452
$ws.Range("M16").Value = 16.02593876639145
$ws.Range("B17").Value = 14.2329777839743
$ws.Range("C17").Value = 6.347076775284452
$ws.Range("D17").Value = 7.097987164550582
$ws.Range("F17").Value = 43.64516049218765
$ws.Range("G17").Value = 52.74204711591537
$ws.Range("H17").Value = 20.25378685759218
$ws.Range("I17").Value = 32.41935286130569
$ws.Range("J17").Value = 10.74888921146227
$ws.Range("K17").Value = 11.69052504887505
$ws.Range("L17").Value = 11.76516930010977
$ws.Range("M17").Value = 16.01613501664386
$ws.Range("B18").Value = 14.19121603931179
$ws.Range("C18").Value = 6.322630007998077
$ws.Range("D18").Value = 7.096141024281935
$ws.Range("F18").Value = 43.65048663741195
$ws.Range("G18").Value = 52.74048900002689
$ws.Range("H18").Value = 20.26190925674298
$ws.Range("I18").Value = 32.43287971857912
$ws.Range("J18").Value = 10.75269748842164
$ws.Range("K18").Value = 11.66051742743641
$ws.Range("L18").Value = 11.766680186642
$ws.Range("M18").Value = 16.01071162961728
$ws.Range("B19").Value = 14.17709441078122
$ws.Range("C19").Value = 6.31431586062354
$ws.Range("D19").Value = 7.095528773083773
$ws.Range("F19").Value = 43.65242274327907
$ws.Range("G19").Value = 52.74014969062779
$ws.Range("H19").Value = 20.26470673152453
$ws.Range("I19").Value = 32.43754226898373
$ws.Range("J19").Value = 10.75399679186042
$ws.Range("K19").Value = 11.65037682228984
$ws.Range("L19").Value = 11.76721432448603
$ws.Range("M19").Value = 16.00891248849209
$ws.Range("B20").Value = 14.24071526664292
$ws.Range("C20").Value = 6.351583714852558
$ws.Range("D20").Value = 7.098334923305563
$ws.Range("F20").Value = 43.64423780220738
$ws.Range("G20").Value = 52.74242489389641
$ws.Range("H20").Value = 20.25230608538874
$ws.Range("I20").Value = 32.41688857867065
$ws.Range("J20").Value = 10.74818907821735
$ws.Range("K20").Value = 11.69608781628978
$ws.Range("L20").Value = 11.76490038790556
$ws.Range("M20").Value = 16.01715637019979
$ws.Range("B21").Value = 14.45554781003375
$ws.Range("C21").Value = 6.474182093289993
$ws.Range("D21").Value = 7.108648953156645
$ws.Range("F21").Value = 43.62597073716699
$ws.Range("G21").Value = 52.7631573876028
$ws.Range("H21").Value = 20.21393736109019
$ws.Range("I21").Value = 32.35324340200228
$ws.Range("J21").Value = 10.72936208270233
$ws.Range("K21").Value = 11.85089545989189
$ws.Range("L21").Value = 11.75869318113958
$ws.Range("M21").Value = 16.04742017847026
$ws.Range("B22").Value = 14.59653313973885
$ws.Range("C22").Value = 6.552300236166733
$ws.Range("D22").Value = 7.116042588434858
$ws.Range("F22").Value = 43.62090916738261
$ws.Range("G22").Value = 52.78645029037089
$ws.Range("H22").Value = 20.19131796441822
$ws.Range("I22").Value = 32.31593533549045
$ws.Range("J22").Value = 10.71756704807164
$ws.Range("K22").Value = 11.95283293170274
$ws.Range("L22").Value = 11.75580052094421
$ws.Range("M22").Value = 16.06909062294486
$ws.Range("B23").Value = 14.5212476888013
$ws.Range("C23").Value = 6.510794291468694
$ws.Range("D23").Value = 7.112037214223052
$ws.Range("F23").Value = 43.62298176273442
$ws.Range("G23").Value = 52.77312789366005
$ws.Range("H23").Value = 20.20316554697104
$ws.Range("I23").Value = 32.33545482335125
$ws.Range("J23").Value = 10.72381578165667
$ws.Range("K23").Value = 11.89836697715275
$ws.Range("L23").Value = 11.75723813803214
$ws.Range("M23").Value = 16.05735301977078
$ws.Range("B24").Value = 14.23721690069184
$ws.Range("C24").Value = 6.349546836265094
$ws.Range("D24").Value = 7.098177472441826
$ws.Range("F24").Value = 43.64465253622109
$ws.Range("G24").Value = 52.74225069071555
$ws.Range("H24").Value = 20.25297467223618
$ws.Range("I24").Value = 32.4180011642067
$ws.Range("J24").Value = 10.7485054240462
$ws.Range("K24").Value = 11.69357259267889
$ws.Range("L24").Value = 11.76502155205472
$ws.Range("M24").Value = 16.01669395273498
$ws.Range("B25").Value = 13.93477289272715
$ws.Range("C25").Value = 6.167375543599278
$ws.Range("D25").Value = 7.086076997330561
$ws.Range("F25").Value = 43.6975761876475
$ws.Range("G25").Value = 52.75083499678426
$ws.Range("H25").Value = 20.31724653162081
$ws.Range("I25").Value = 32.5254265389506
$ws.Range("J25").Value = 10.77733697667344
$ws.Range("K25").Value = 11.47691876126244
$ws.Range("L25").Value = 11.77844063064837
$ws.Range("M25").Value = 15.98108806079533
